$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 of the BOM: the ADC input RC resistor changed from
# 50 Ohms 0.1% (RT0402BRE0750RL) to 200 Ohms 0.5% (RT0402DRE07200RL).
$ws.Range("A12").Value = "RES SMD 200 OHM 0.5% 1/16W 0402"
$ws.Range("B12").Value = "200 Ohms ±0.5% 0.063W, 1/16W Chip Resistor 0402 (1005 Metric) Thin Film"

# New manufacturer part number column value for this row, wrapped like a
# multi-line cell.
$ws.Range("E12").Value = "RT0402DRE07200RL"
$ws.Range("E12").WrapText = $true

# Update the Digi-Key hyperlink (and its displayed text) to point at the
# new part, and give it the same Hyperlink style used elsewhere.
$hl = $ws.Hyperlinks.Item(16)
$hl.Address = "https://www.digikey.com/en/products/detail/yageo/RT0402DRE07200RL/1071388"
$hl.TextToDisplay = "https://www.digikey.com/en/products/detail/yageo/RT0402DRE07200RL/1071388"
$ws.Range("H12").Style = "Hyperlink"
